$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Narration"
$ws.Range("G2").Value = "Here we go"
$ws.Range("G3").Value = "Here also"
$ws.Range("G4").Value = "Finally here"

$ws.Range("G4").Select()
